# Create prediction files for 5 years.
#
# The workbook originally tracked a separate "government" account in its own
# column (column B, header "government", data value "GOV").  The new layout
# folds the government account into the household/agent list in column A
# (after the HOU_* rows) instead of keeping a dedicated column for it.
#
# Net effect:
#   1. Column B ("government" / "GOV" / blanks) is removed entirely; the
#      former columns C:F (factors, inc_taxes, goods_activities,
#      goods_commodities and all their rows) shift one column to the left
#      to become B:E.
#   2. The "GOV" label is re-added as a new row entry in column A,
#      immediately after the HOU_95-100 row (row 7), matching the other
#      plain (unstyled) label cells in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "government" column (column B) - this shifts C:F left to B:E
#    and automatically drops the now-unused "government"/"GOV" shared
#    strings.
$ws.Columns("B:B").Delete()

# 2. Re-insert the government row label in column A, right after the last
#    household row (HOU_95-100 is A6), as a plain (unstyled) text cell.
$ws.Range("A7").Value = "GOV"

# Reset the selection so it no longer references the old (now invalid)
# F2:F22 range from before the column was removed.
$ws.Range("A1").Select() | Out-Null
